$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation needs to be inserted as row 22 (Fecha 2022-07-22),
# pushing every following record down by one row (old row 22 -> new row 23,
# ..., old row 91 -> new row 92). Excel's native row insert takes care of
# shifting all the existing values/styles down for us.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new data point. It mirrors the
# record that used to sit at row 22 (same Volumen/Origen), but with its own
# Fecha and Precio values.
$ws.Range("A22").Value = 4
$ws.Range("B22").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C22").Value = "Los Lagos"
$ws.Range("D22").Value = 44764
$ws.Range("E22").Value = 10
$ws.Range("F22").Value = 100112026
$ws.Range("G22").Value = "Haba"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 80
$ws.Range("K22").Value = 19000
$ws.Range("L22").Value = 19000
$ws.Range("M22").Value = 19000
$ws.Range("N22").Value = "$/saco 25 kilos"
$ws.Range("O22").Value = "Provincia de Limarí"
$ws.Range("P22").Value = 760
$ws.Range("Q22").Value = 25
$ws.Range("R22").Value = "Hortaliza"
